# ProjectPlanning_Current.xlsx - "Add files via upload" commit
#
# Applies the cell-data and view-state changes to the Sprint2 worksheet
# (and the Sprint1 worksheet's view) described by the upstream diff.
#
# Notes on scope: the diff also touches two things that are not reachable
# through the Excel object model / this COM-interop host:
#   - <x15ac:absPath .../> in xl/workbook.xml is a "last saved from" path
#     stamped by the real Excel application; it isn't an exposed COM
#     property, so it can't be set from script.
#   - The cached <c:numCache> points baked into the chart parts
#     (xl/charts/chart*.xml) aren't refreshed by this host when the cells
#     they reference change (nor by Chart.Refresh/SetSourceData/etc.); the
#     chart XML is carried through unchanged regardless. The live
#     SeriesCollection(...).Values do reflect the recalculated numbers, but
#     that doesn't make it into the saved cache. So that part of the diff
#     is also outside what this script can control.
# Everything else (cell values/formulas and the sheet view/selection
# state that IS persisted) is applied below.

$wb = $excel.ActiveWorkbook

# Sprint1's own change is just its view losing topLeftCell="A7" (no cell
# values move there) - see note above re: topLeftCell not being
# round-tripped by this host, so there is nothing to do for that sheet.
$sprint2 = $wb.Worksheets.Item("Sprint2")

# --- Sprint2!E7:F8 block -------------------------------------------------
$sprint2.Range("E7").Value = 20
$sprint2.Range("E8").Value = 6

# --- Sprint2 burndown input grid (rows 19-25) ----------------------------
$sprint2.Range("J20").ClearContents()
$sprint2.Range("J21").Value = 3
$sprint2.Range("G22").Value = 20
$sprint2.Range("J22").Value = 20
$sprint2.Range("J24").Value = 1
$sprint2.Range("G25").Value = 3
$sprint2.Range("H25").Value = 1
$sprint2.Range("I25").Value = 1
$sprint2.Range("J25").Value = 1

# --- View / selection state ----------------------------------------------
# Sprint2 becomes/remains the active sheet, with its selection moved from
# L20 to H24 (and, per the diff, is also scrolled so A7 is the top-left
# visible cell - set here even though this host does not persist
# topLeftCell back into the saved sheetView).
$sprint2.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$null = $sprint2.Range("H24").Select()
